$wb = $excel.ActiveWorkbook

# Sheet 1: "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 50
$ws1.Range("F5").Value = 186
$ws1.Range("F6").Value = 1085
$ws1.Range("F7").Value = 1048
$ws1.Range("F8").Value = 8195
$ws1.Range("F9").Value = 138
$ws1.Range("F10").Value = 209
$ws1.Range("F11").Value = 6901
$ws1.Range("F12").Value = 173
$ws1.Range("F13").Value = 301
$ws1.Range("F14").Value = 5006
$ws1.Range("F16").Value = 5456
$ws1.Range("F17").Value = 1074
$ws1.Range("F18").Value = 333
$ws1.Range("F19").Value = 343
$ws1.Range("F20").Value = 467
$ws1.Range("F26").Value = 9199
$ws1.Range("F28").Value = 1672
$ws1.Range("F29").Value = 819
$ws1.Range("F30").Value = 39
$ws1.Range("F32").Value = 1869
$ws1.Range("F34").Value = 79
$ws1.Range("F37").Value = 1884
$ws1.Range("F40").Value = 4805
$ws1.Range("F42").Value = 1163
$ws1.Range("F43").Value = 77
$ws1.Range("F44").Value = 148
$ws1.Range("F45").Value = 74
$ws1.Range("F46").Value = 37
$ws1.Range("F47").Value = 918
$ws1.Range("F48").Value = 1260
$ws1.Range("F49").Value = 41
$ws1.Range("F50").Value = 64

# Sheet 2: "演出" (Performance)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F3").Value = 40
$ws2.Range("F6").Value = 24
$ws2.Range("F17").Value = 891

# Sheet 4: "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value = 50
$ws4.Range("F6").Value = 186
$ws4.Range("F7").Value = 40
$ws4.Range("F8").Value = 1085
$ws4.Range("F9").Value = 1048
$ws4.Range("F10").Value = 8195
$ws4.Range("F11").Value = 138
$ws4.Range("F12").Value = 209
$ws4.Range("F13").Value = 6901
$ws4.Range("F14").Value = 173
$ws4.Range("F15").Value = 301
$ws4.Range("F17").Value = 5006
$ws4.Range("F19").Value = 5456
$ws4.Range("F20").Value = 1074
$ws4.Range("F21").Value = 333
$ws4.Range("F22").Value = 343
$ws4.Range("F23").Value = 467
$ws4.Range("F27").Value = 9199
$ws4.Range("F29").Value = 1672
$ws4.Range("F30").Value = 819
$ws4.Range("F31").Value = 39
$ws4.Range("F33").Value = 1869
$ws4.Range("F35").Value = 79
$ws4.Range("F38").Value = 1884
$ws4.Range("F41").Value = 4805
$ws4.Range("F43").Value = 1163
$ws4.Range("F44").Value = 77
$ws4.Range("F45").Value = 148
$ws4.Range("F46").Value = 74
$ws4.Range("F47").Value = 918
$ws4.Range("F48").Value = 1260
$ws4.Range("F49").Value = 41
$ws4.Range("F50").Value = 64
